# Update F2-F2rl2 LR-pairs sheet with re-run TPM numbers.
# - Rows that used "Resolving-Mac" as a sending cluster paired with every target
#   cluster previously lived at rows 14-17; the new TPM run folds those pairs
#   into rows 11-13 instead (Target cluster "ECs" is dropped entirely), so the
#   last four rows are deleted and every remaining data row gets new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("14:17").Delete()

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F2"
$ws.Range("C2").Value = "F2rl2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.754521
$ws.Range("H2").Value = 2.263563
$ws.Range("I2").Value = 0.2768403531129761
$ws.Range("J2").Value = 0.2768403531129761
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3700743333333333
$ws.Range("N2").Value = 1.110223
$ws.Range("O2").Value = 0.1108002562864021
$ws.Range("P2").Value = 0.1108002562864021
$ws.Range("Q2").Value = 0.279228856061
$ws.Range("R2").Value = 2.513059704549
$ws.Range("S2").Value = 0.0306739820753358
$ws.Range("T2").Value = 0.0306739820753358

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F2"
$ws.Range("C3").Value = "F2rl2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.754521
$ws.Range("H3").Value = 2.263563
$ws.Range("I3").Value = 0.2768403531129761
$ws.Range("J3").Value = 0.2768403531129761
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8984456666666668
$ws.Range("N3").Value = 2.695337
$ws.Range("O3").Value = 0.2689946347519571
$ws.Range("P3").Value = 0.2689946347519571
$ws.Range("Q3").Value = 0.6778961228590001
$ws.Range("R3").Value = 6.101065105731001
$ws.Range("S3").Value = 0.07446856967022784
$ws.Range("T3").Value = 0.07446856967022784

# Row 4: ECs -> Resolving-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "F2"
$ws.Range("C4").Value = "F2rl2"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.754521
$ws.Range("H4").Value = 2.263563
$ws.Range("I4").Value = 0.2768403531129761
$ws.Range("J4").Value = 0.2768403531129761
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.071493333333333
$ws.Range("N4").Value = 6.21448
$ws.Range("O4").Value = 0.6202051089616408
$ws.Range("P4").Value = 0.6202051089616408
$ws.Range("Q4").Value = 1.56298522136
$ws.Range("R4").Value = 14.06686699224
$ws.Range("S4").Value = 0.1716978013674125
$ws.Range("T4").Value = 0.1716978013674125

# Row 5: FAPs -> FAPs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F2"
$ws.Range("C5").Value = "F2rl2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9731926666666667
$ws.Range("H5").Value = 2.919578
$ws.Range("I5").Value = 0.3570728998754956
$ws.Range("J5").Value = 0.3570728998754956
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3700743333333333
$ws.Range("N5").Value = 1.110223
$ws.Range("O5").Value = 0.1108002562864021
$ws.Range("P5").Value = 0.1108002562864021
$ws.Range("Q5").Value = 0.3601536273215556
$ws.Range("R5").Value = 3.241382645894
$ws.Range("S5").Value = 0.0395637688191337
$ws.Range("T5").Value = 0.0395637688191337

# Row 6: FAPs -> MuSCs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "F2"
$ws.Range("C6").Value = "F2rl2"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9731926666666667
$ws.Range("H6").Value = 2.919578
$ws.Range("I6").Value = 0.3570728998754956
$ws.Range("J6").Value = 0.3570728998754956
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8984456666666668
$ws.Range("N6").Value = 2.695337
$ws.Range("O6").Value = 0.2689946347519571
$ws.Range("P6").Value = 0.2689946347519571
$ws.Range("Q6").Value = 0.8743607341984446
$ws.Range("R6").Value = 7.869246607786001
$ws.Range("S6").Value = 0.09605069428183109
$ws.Range("T6").Value = 0.09605069428183109

# Row 7: FAPs -> Resolving-Mac
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "F2"
$ws.Range("C7").Value = "F2rl2"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9731926666666667
$ws.Range("H7").Value = 2.919578
$ws.Range("I7").Value = 0.3570728998754956
$ws.Range("J7").Value = 0.3570728998754956
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.071493333333333
$ws.Range("N7").Value = 6.21448
$ws.Range("O7").Value = 0.6202051089616408
$ws.Range("P7").Value = 0.6202051089616408
$ws.Range("Q7").Value = 2.015962121048889
$ws.Range("R7").Value = 18.14365908944
$ws.Range("S7").Value = 0.2214584367745308
$ws.Range("T7").Value = 0.2214584367745308

# Row 8: MuSCs -> FAPs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "F2"
$ws.Range("C8").Value = "F2rl2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.782441
$ws.Range("H8").Value = 2.347323
$ws.Range("I8").Value = 0.2870844452706686
$ws.Range("J8").Value = 0.2870844452706686
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.3700743333333333
$ws.Range("N8").Value = 1.110223
$ws.Range("O8").Value = 0.1108002562864021
$ws.Range("P8").Value = 0.1108002562864021
$ws.Range("Q8").Value = 0.2895613314476667
$ws.Range("R8").Value = 2.606051983029
$ws.Range("S8").Value = 0.03180903011182964
$ws.Range("T8").Value = 0.03180903011182964

# Row 9: MuSCs -> MuSCs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "F2"
$ws.Range("C9").Value = "F2rl2"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.782441
$ws.Range("H9").Value = 2.347323
$ws.Range("I9").Value = 0.2870844452706686
$ws.Range("J9").Value = 0.2870844452706686
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8984456666666668
$ws.Range("N9").Value = 2.695337
$ws.Range("O9").Value = 0.2689946347519571
$ws.Range("P9").Value = 0.2689946347519571
$ws.Range("Q9").Value = 0.7029807258723335
$ws.Range("R9").Value = 6.326826532851001
$ws.Range("S9").Value = 0.07722417549855172
$ws.Range("T9").Value = 0.07722417549855173

# Row 10: MuSCs -> Resolving-Mac
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "F2"
$ws.Range("C10").Value = "F2rl2"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.782441
$ws.Range("H10").Value = 2.347323
$ws.Range("I10").Value = 0.2870844452706686
$ws.Range("J10").Value = 0.2870844452706686
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.071493333333333
$ws.Range("N10").Value = 6.21448
$ws.Range("O10").Value = 0.6202051089616408
$ws.Range("P10").Value = 0.6202051089616408
$ws.Range("Q10").Value = 1.620821315226667
$ws.Range("R10").Value = 14.58739183704
$ws.Range("S10").Value = 0.1780512396602872
$ws.Range("T10").Value = 0.1780512396602872

# Row 11: Resolving-Mac -> FAPs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "F2"
$ws.Range("C11").Value = "F2rl2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2153186666666667
$ws.Range("H11").Value = 0.645956
$ws.Range("I11").Value = 0.07900230174085969
$ws.Range("J11").Value = 0.07900230174085969
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.3700743333333333
$ws.Range("N11").Value = 1.110223
$ws.Range("O11").Value = 0.1108002562864021
$ws.Range("P11").Value = 0.1108002562864021
$ws.Range("Q11").Value = 0.07968391202088888
$ws.Range("R11").Value = 0.7171552081879999
$ws.Range("S11").Value = 0.00875347528010292
$ws.Range("T11").Value = 0.00875347528010292

# Row 12: Resolving-Mac -> MuSCs
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "F2"
$ws.Range("C12").Value = "F2rl2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.2153186666666667
$ws.Range("H12").Value = 0.645956
$ws.Range("I12").Value = 0.07900230174085969
$ws.Range("J12").Value = 0.07900230174085969
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.8984456666666668
$ws.Range("N12").Value = 2.695337
$ws.Range("O12").Value = 0.2689946347519571
$ws.Range("P12").Value = 0.2689946347519571
$ws.Range("Q12").Value = 0.1934521230191111
$ws.Range("R12").Value = 1.741069107172
$ws.Range("S12").Value = 0.02125119530134646
$ws.Range("T12").Value = 0.02125119530134646

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "F2"
$ws.Range("C13").Value = "F2rl2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.2153186666666667
$ws.Range("H13").Value = 0.645956
$ws.Range("I13").Value = 0.07900230174085969
$ws.Range("J13").Value = 0.07900230174085969
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.071493333333333
$ws.Range("N13").Value = 6.21448
$ws.Range("O13").Value = 0.6202051089616408
$ws.Range("P13").Value = 0.6202051089616408
$ws.Range("Q13").Value = 0.4460311825422222
$ws.Range("R13").Value = 4.01428064288
$ws.Range("S13").Value = 0.04899763115941031
$ws.Range("T13").Value = 0.04899763115941031
